# Adds the new match row (row 4) to the active worksheet, replicating
# the structure/values of the existing rows (Id, Date, Time, League,
# Home, Away team, and all odds columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "rFf4UJrf"
$ws.Range("B4").Value = "30/10/2024"
$ws.Range("C4").Value = "12:00"
$ws.Range("D4").Value = "SLOVAKIA - NIKE LIGA"
$ws.Range("E4").Value = "Skalica"
$ws.Range("F4").Value = "Slovan Bratislava"
$ws.Range("G4").Value = 5.1
$ws.Range("H4").Value = 4.35
$ws.Range("I4").Value = 1.55
$ws.Range("J4").Value = 4.9
$ws.Range("K4").Value = 2.45
$ws.Range("L4").Value = 2.02
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 9.75
$ws.Range("O4").Value = 1.19
$ws.Range("P4").Value = 4.4
$ws.Range("Q4").Value = 1.57
$ws.Range("R4").Value = 2.35
$ws.Range("S4").Value = 1.29
$ws.Range("T4").Value = 3.35
$ws.Range("U4").Value = 1.65
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 16
$ws.Range("X4").Value = 35
$ws.Range("Y4").Value = 17.5
$ws.Range("Z4").Value = 100
$ws.Range("AA4").Value = 50
$ws.Range("AB4").Value = 45
$ws.Range("AC4").Value = 9.75
$ws.Range("AD4").Value = 9
$ws.Range("AE4").Value = 16
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 8.25
$ws.Range("AI4").Value = 9
$ws.Range("AJ4").Value = 8.75
$ws.Range("AK4").Value = 12.5
$ws.Range("AL4").Value = 12
$ws.Range("AM4").Value = 22
$ws.Range("AN4").Value = 6.9
$ws.Range("AO4").Value = 27
$ws.Range("AP4").Value = 28
$ws.Range("AQ4").Value = 150
$ws.Range("AR4").Value = 150
$ws.Range("AS4").Value = 300
$ws.Range("AT4").Value = 3.35
$ws.Range("AU4").Value = 7.3
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 3.6
$ws.Range("AX4").Value = 7.1
$ws.Range("AY4").Value = 14.5
$ws.Range("AZ4").Value = 20
$ws.Range("BA4").Value = 40
$ws.Range("BB4").Value = 150
$ws.Range("BC4").Value = 51
$ws.Range("BD4").Value = 51